# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data table with the latest figures and
# update the "last updated" timestamp. A couple of countries
# (Emiratos Arabes Unidos / Dinamarca and Serbia / Chequia) swap rank
# order because their totals changed relative to each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Abril de 2020 a las 15:22"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 819321
$ws.Cells.Item(4, 3).Value = 577
$ws.Cells.Item(4, 5).Value = 690993
$ws.Cells.Item(4, 7).Value = 37
$ws.Cells.Item(4, 8).Value = 45355

# --- Alemania (row 8) ---
$ws.Cells.Item(8, 2).Value = 148766
$ws.Cells.Item(8, 3).Value = 313
$ws.Cells.Item(8, 5).Value = 44264

# --- Portugal (row 19) ---
$ws.Cells.Item(19, 6).Value = 207

# --- Arabia Saudita (row 26) ---
$ws.Cells.Item(26, 6).Value = 82

# --- Dinamarca / Emiratos Arabes Unidos swap (rows 36-37) ---
# Emiratos Arabes Unidos overtakes Dinamarca with new totals.
$ws.Cells.Item(36, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(36, 2).Value = 8238
$ws.Cells.Item(36, 3).Value = 483
$ws.Cells.Item(36, 4).Value = 1546
$ws.Cells.Item(36, 5).Value = 6640
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 6
$ws.Cells.Item(36, 8).Value = 52

$ws.Cells.Item(37, 1).Value = "Dinamarca"
$ws.Cells.Item(37, 2).Value = 7912
$ws.Cells.Item(37, 3).Value = 217
$ws.Cells.Item(37, 4).Value = 5087
$ws.Cells.Item(37, 5).Value = 2441
$ws.Cells.Item(37, 6).Value = 80
$ws.Cells.Item(37, 7).Value = 14
$ws.Cells.Item(37, 8).Value = 384

# --- Chequia / Serbia swap (rows 42-43) ---
# Serbia overtakes Chequia with new totals.
$ws.Cells.Item(42, 1).Value = "Serbia"
$ws.Cells.Item(42, 2).Value = 7114
$ws.Cells.Item(42, 3).Value = 224
$ws.Cells.Item(42, 4).Value = 1025
$ws.Cells.Item(42, 5).Value = 5955
$ws.Cells.Item(42, 6).Value = 101
$ws.Cells.Item(42, 7).Value = 4
$ws.Cells.Item(42, 8).Value = 134

$ws.Cells.Item(43, 1).Value = "Chequia"
$ws.Cells.Item(43, 2).Value = 7041
$ws.Cells.Item(43, 3).Value = 8
$ws.Cells.Item(43, 4).Value = 1800
$ws.Cells.Item(43, 5).Value = 5037
$ws.Cells.Item(43, 6).Value = 80
$ws.Cells.Item(43, 7).Value = 3
$ws.Cells.Item(43, 8).Value = 204

# --- Republica de Yibuti (row 86) ---
$ws.Cells.Item(86, 2).Value = 974
$ws.Cells.Item(86, 3).Value = 29
$ws.Cells.Item(86, 4).Value = 183
$ws.Cells.Item(86, 5).Value = 789

# --- San Marino (row 104) ---
$ws.Cells.Item(104, 2).Value = 488
$ws.Cells.Item(104, 3).Value = 12
$ws.Cells.Item(104, 5).Value = 386

# --- Kenia (row 119) ---
$ws.Cells.Item(119, 2).Value = 303
$ws.Cells.Item(119, 3).Value = 7
$ws.Cells.Item(119, 5).Value = 215
